$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: "A form to express interest and gather some information
# about the user." -> "A form for anonymous feed back to the teachers."
# split across three runs: "A form " | "for anonymous feed back to the
# teachers" | "."
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute(
    "A form to express interest and gather some information about the user.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A form for anonymous feed back to the teachers.", 2)

# locate the paragraph that now holds the replaced sentence
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "A form for anonymous feed back to the teachers.*") {
        $p1 = $p
        break
    }
}

$pStart = $p1.Range.Start
$pEnd = $p1.Range.End
$textEnd = $pEnd - 1   # exclude the paragraph mark

# Split off the trailing "." into its own run (do the right-most split
# first so earlier offsets stay valid).
$splitB = $d.Range($textEnd - 1, $textEnd)
$splitB.Font.Bold = 1
$splitB.Font.Bold = 0

# Split "A form " from "for anonymous feed back to the teachers"
$splitA = $d.Range($pStart, $pStart + 7)
$splitA.Font.Bold = 1
$splitA.Font.Bold = 0

# ---------------------------------------------------------------------
# Change 2: "I do not see a form for gathering informatio" + "n about a
# user if interested" (two runs separated by the _GoBack bookmark)
# -> "I do not see a form for " | "feed" + bookmark + "back"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "I do not see a form for gathering information about a user if interested*") {
        $p2 = $p
        break
    }
}

$p2Start = $p2.Range.Start

# First run, before the bookmark: "I do not see a form for gathering informatio"
$run1 = $d.Range($p2Start, $p2Start + 44)
$run1.Text = "I do not see a form for feed"

$afterRun1 = $p2Start + ("I do not see a form for feed").Length

# Second run, after the bookmark, up to (but excluding) the paragraph mark
$p2End = $p2.Range.End
$run2 = $d.Range($afterRun1, $p2End - 1)
$run2.Text = "back"

# Split "I do not see a form for feed" into "I do not see a form for "
# and "feed" (two separate runs).
$splitPos = $p2Start + ("I do not see a form for ").Length
$splitC = $d.Range($splitPos, $afterRun1)
$splitC.Font.Bold = 1
$splitC.Font.Bold = 0
